$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("B17").Value = 6851963
$ws.Range("E17").Value = "FC Blau Weiss Linz"
$ws.Range("F17").Value = "Hartberg"
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 3
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = "D"
$ws.Range("L17").Value = 2.7
$ws.Range("M17").Value = 3.4
$ws.Range("N17").Value = 2.55
$ws.Range("O17").Value = 2.25
$ws.Range("P17").Value = 3.5
$ws.Range("Q17").Value = 3.1
$ws.Range("R17").Value = -0.25
$ws.Range("S17").Value = 1.975
$ws.Range("T17").Value = 1.875
$ws.Range("U17").Value = 2.75
$ws.Range("V17").Value = 1.9
$ws.Range("W17").Value = 1.95
$ws.Range("X17").Value = -1
$ws.Range("Y17").Value = 2.5
$ws.Range("Z17").Value = -1
$ws.Range("AA17").Value = -0.5
$ws.Range("AB17").Value = 0.4375
$ws.Range("AC17").Value = 0.8999999999999999
$ws.Range("AD17").Value = -1

# Row 18
$ws.Range("B18").Value = 6847032
$ws.Range("E18").Value = "Austria Lustenau"
$ws.Range("F18").Value = "FK Austria Vienna"
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 2
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = "A"
$ws.Range("L18").Value = 4.2
$ws.Range("M18").Value = 4
$ws.Range("N18").Value = 1.75
$ws.Range("O18").Value = 3.75
$ws.Range("P18").Value = 3.8
$ws.Range("Q18").Value = 1.909
$ws.Range("R18").Value = 0.5
$ws.Range("S18").Value = 1.975
$ws.Range("T18").Value = 1.875
$ws.Range("U18").Value = 3
$ws.Range("V18").Value = 1.95
$ws.Range("W18").Value = 1.9
$ws.Range("X18").Value = -1
$ws.Range("Y18").Value = -1
$ws.Range("Z18").Value = 0.909
$ws.Range("AA18").Value = -1
$ws.Range("AB18").Value = 0.875
$ws.Range("AC18").Value = -1
$ws.Range("AD18").Value = 0.8999999999999999

# Row 44
$ws.Range("B44").Value = 6847049
$ws.Range("E44").Value = "Austria Klagenfurt"
$ws.Range("F44").Value = "LASK Linz"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 3
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 3
$ws.Range("K44").Value = "A"
$ws.Range("L44").Value = 3.4
$ws.Range("M44").Value = 3.6
$ws.Range("N44").Value = 1.909
$ws.Range("O44").Value = 3.3
$ws.Range("P44").Value = 3.4
$ws.Range("Q44").Value = 2.2
$ws.Range("R44").Value = 0.25
$ws.Range("S44").Value = 1.95
$ws.Range("T44").Value = 1.9
$ws.Range("U44").Value = 2.5
$ws.Range("V44").Value = 1.975
$ws.Range("W44").Value = 1.875
$ws.Range("X44").Value = -1
$ws.Range("Y44").Value = -1
$ws.Range("Z44").Value = 1.2
$ws.Range("AA44").Value = -1
$ws.Range("AB44").Value = 0.8999999999999999
$ws.Range("AC44").Value = 0.9750000000000001
$ws.Range("AD44").Value = -1

# Row 45
$ws.Range("B45").Value = 6851958
$ws.Range("E45").Value = "FC Blau Weiss Linz"
$ws.Range("F45").Value = "SCR Altach"
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 1
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = 1
$ws.Range("K45").Value = "D"
$ws.Range("L45").Value = 2.3
$ws.Range("M45").Value = 3.5
$ws.Range("N45").Value = 2.7
$ws.Range("O45").Value = 2.375
$ws.Range("P45").Value = 3.4
$ws.Range("Q45").Value = 3
$ws.Range("R45").Value = -0.25
$ws.Range("S45").Value = 2.025
$ws.Range("T45").Value = 1.825
$ws.Range("U45").Value = 2.5
$ws.Range("V45").Value = 1.975
$ws.Range("W45").Value = 1.875
$ws.Range("X45").Value = -1
$ws.Range("Y45").Value = 2.4
$ws.Range("Z45").Value = -1
$ws.Range("AA45").Value = -0.5
$ws.Range("AB45").Value = 0.4125
$ws.Range("AC45").Value = -1
$ws.Range("AD45").Value = 0.875

# Row 75
$ws.Range("B75").Value = 6851953
$ws.Range("E75").Value = "FC Blau Weiss Linz"
$ws.Range("F75").Value = "Wolfsberger AC"
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = "H"
$ws.Range("L75").Value = 2.7
$ws.Range("M75").Value = 3.5
$ws.Range("N75").Value = 2.3
$ws.Range("O75").Value = 2.625
$ws.Range("P75").Value = 3.5
$ws.Range("Q75").Value = 2.55
$ws.Range("R75").Value = 0
$ws.Range("S75").Value = 1.95
$ws.Range("T75").Value = 1.9
$ws.Range("U75").Value = 2.5
$ws.Range("V75").Value = 1.925
$ws.Range("W75").Value = 1.925
$ws.Range("X75").Value = 1.625
$ws.Range("Y75").Value = -1
$ws.Range("Z75").Value = -1
$ws.Range("AA75").Value = 0.95
$ws.Range("AB75").Value = -1
$ws.Range("AC75").Value = -1
$ws.Range("AD75").Value = 0.925

# Row 76
$ws.Range("B76").Value = 6847071
$ws.Range("E76").Value = "Austria Klagenfurt"
$ws.Range("F76").Value = "WSG Swarovski Tirol"
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = "H"
$ws.Range("L76").Value = 1.727
$ws.Range("M76").Value = 3.75
$ws.Range("N76").Value = 4
$ws.Range("O76").Value = 1.75
$ws.Range("P76").Value = 3.75
$ws.Range("Q76").Value = 4.75
$ws.Range("R76").Value = -0.75
$ws.Range("S76").Value = 1.975
$ws.Range("T76").Value = 1.875
$ws.Range("U76").Value = 2.75
$ws.Range("V76").Value = 2
$ws.Range("W76").Value = 1.85
$ws.Range("X76").Value = 0.75
$ws.Range("Y76").Value = -1
$ws.Range("Z76").Value = -1
$ws.Range("AA76").Value = 0.4875
$ws.Range("AB76").Value = -0.5
$ws.Range("AC76").Value = -1
$ws.Range("AD76").Value = 0.8500000000000001

# Row 77
$ws.Range("B77").Value = 6847073
$ws.Range("E77").Value = "Austria Lustenau"
$ws.Range("F77").Value = "Hartberg"
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1
$ws.Range("K77").Value = "A"
$ws.Range("L77").Value = 3.4
$ws.Range("M77").Value = 3.5
$ws.Range("N77").Value = 1.95
$ws.Range("O77").Value = 4.5
$ws.Range("P77").Value = 4
$ws.Range("Q77").Value = 1.7
$ws.Range("R77").Value = 0.75
$ws.Range("S77").Value = 1.9
$ws.Range("T77").Value = 1.95
$ws.Range("U77").Value = 2.75
$ws.Range("V77").Value = 1.875
$ws.Range("W77").Value = 1.975
$ws.Range("X77").Value = -1
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = 0.7
$ws.Range("AA77").Value = -1
$ws.Range("AB77").Value = 0.95
$ws.Range("AC77").Value = 0.875
$ws.Range("AD77").Value = -1

# Row 78
$ws.Range("B78").Value = 6847070
$ws.Range("E78").Value = "SK Sturm Graz"
$ws.Range("F78").Value = "FK Austria Vienna"
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 1
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 1
$ws.Range("K78").Value = "A"
$ws.Range("L78").Value = 1.5
$ws.Range("M78").Value = 3.8
$ws.Range("N78").Value = 6
$ws.Range("O78").Value = 1.727
$ws.Range("P78").Value = 3.8
$ws.Range("Q78").Value = 4.75
$ws.Range("R78").Value = -0.75
$ws.Range("S78").Value = 1.95
$ws.Range("T78").Value = 1.9
$ws.Range("U78").Value = 2.75
$ws.Range("V78").Value = 1.875
$ws.Range("W78").Value = 1.975
$ws.Range("X78").Value = -1
$ws.Range("Y78").Value = -1
$ws.Range("Z78").Value = 3.75
$ws.Range("AA78").Value = -1
$ws.Range("AB78").Value = 0.8999999999999999
$ws.Range("AC78").Value = -1
$ws.Range("AD78").Value = 0.9750000000000001

# Row 86
$ws.Range("B86").Value = 6847081
$ws.Range("E86").Value = "WSG Swarovski Tirol"
$ws.Range("F86").Value = "SCR Altach"
$ws.Range("G86").Value = 5
$ws.Range("H86").Value = 1
$ws.Range("I86").Value = 3
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = "H"
$ws.Range("L86").Value = 2.8
$ws.Range("M86").Value = 3.4
$ws.Range("N86").Value = 2.4
$ws.Range("O86").Value = 3.6
$ws.Range("P86").Value = 3.5
$ws.Range("Q86").Value = 2.05
$ws.Range("R86").Value = 0.25
$ws.Range("S86").Value = 2.05
$ws.Range("T86").Value = 1.8
$ws.Range("U86").Value = 2.5
$ws.Range("V86").Value = 1.95
$ws.Range("W86").Value = 1.9
$ws.Range("X86").Value = 2.6
$ws.Range("Y86").Value = -1
$ws.Range("Z86").Value = -1
$ws.Range("AA86").Value = 1.05
$ws.Range("AB86").Value = -1
$ws.Range("AC86").Value = 0.95
$ws.Range("AD86").Value = -1

# Row 87
$ws.Range("B87").Value = 6847080
$ws.Range("E87").Value = "Hartberg"
$ws.Range("F87").Value = "Rapid Vienna"
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = "H"
$ws.Range("L87").Value = 3.2
$ws.Range("M87").Value = 3.5
$ws.Range("N87").Value = 2.15
$ws.Range("O87").Value = 3.5
$ws.Range("P87").Value = 3.75
$ws.Range("Q87").Value = 2
$ws.Range("R87").Value = 0.5
$ws.Range("S87").Value = 1.8
$ws.Range("T87").Value = 2.05
$ws.Range("U87").Value = 3
$ws.Range("V87").Value = 2.025
$ws.Range("W87").Value = 1.825
$ws.Range("X87").Value = 2.5
$ws.Range("Y87").Value = -1
$ws.Range("Z87").Value = -1
$ws.Range("AA87").Value = 0.8
$ws.Range("AB87").Value = -1
$ws.Range("AC87").Value = -1
$ws.Range("AD87").Value = 0.825

# Row 113
$ws.Range("B113").Value = 6847097
$ws.Range("E113").Value = "WSG Swarovski Tirol"
$ws.Range("F113").Value = "Austria Lustenau"
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 2
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = "A"
$ws.Range("L113").Value = 1.909
$ws.Range("M113").Value = 3.75
$ws.Range("N113").Value = 3.6
$ws.Range("O113").Value = 1.727
$ws.Range("P113").Value = 3.75
$ws.Range("Q113").Value = 4.75
$ws.Range("R113").Value = -0.75
$ws.Range("S113").Value = 1.925
$ws.Range("T113").Value = 1.925
$ws.Range("U113").Value = 2.75
$ws.Range("V113").Value = 2
$ws.Range("W113").Value = 1.85
$ws.Range("X113").Value = -1
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = 3.75
$ws.Range("AA113").Value = -1
$ws.Range("AB113").Value = 0.925
$ws.Range("AC113").Value = -1
$ws.Range("AD113").Value = 0.8500000000000001

# Row 114
$ws.Range("B114").Value = 6847095
$ws.Range("E114").Value = "Wolfsberger AC"
$ws.Range("F114").Value = "Rapid Vienna"
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 2
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 1
$ws.Range("K114").Value = "A"
$ws.Range("L114").Value = 3.4
$ws.Range("M114").Value = 3.75
$ws.Range("N114").Value = 2
$ws.Range("O114").Value = 3.5
$ws.Range("P114").Value = 3.6
$ws.Range("Q114").Value = 2.05
$ws.Range("R114").Value = 0.5
$ws.Range("S114").Value = 1.825
$ws.Range("T114").Value = 2.025
$ws.Range("U114").Value = 2.75
$ws.Range("V114").Value = 1.925
$ws.Range("W114").Value = 1.925
$ws.Range("X114").Value = -1
$ws.Range("Y114").Value = -1
$ws.Range("Z114").Value = 1.05
$ws.Range("AA114").Value = -1
$ws.Range("AB114").Value = 1.025
$ws.Range("AC114").Value = -1
$ws.Range("AD114").Value = 0.925

# Row 116
$ws.Range("B116").Value = 6847098
$ws.Range("E116").Value = "FK Austria Vienna"
$ws.Range("F116").Value = "SCR Altach"
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 1
$ws.Range("I116").Value = 2
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = "H"
$ws.Range("L116").Value = 1.8
$ws.Range("M116").Value = 3.6
$ws.Range("N116").Value = 4.2
$ws.Range("O116").Value = 1.8
$ws.Range("P116").Value = 3.6
$ws.Range("Q116").Value = 4.5
$ws.Range("R116").Value = -0.5
$ws.Range("S116").Value = 1.825
$ws.Range("T116").Value = 2.025
$ws.Range("U116").Value = 2.5
$ws.Range("V116").Value = 1.975
$ws.Range("W116").Value = 1.875
$ws.Range("X116").Value = 0.8
$ws.Range("Y116").Value = -1
$ws.Range("Z116").Value = -1
$ws.Range("AA116").Value = 0.825
$ws.Range("AB116").Value = -1
$ws.Range("AC116").Value = 0.9750000000000001
$ws.Range("AD116").Value = -1

# Row 117
$ws.Range("B117").Value = 6851942
$ws.Range("E117").Value = "FC Blau Weiss Linz"
$ws.Range("F117").Value = "FC Salzburg"
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 1
$ws.Range("I117").Value = 1
$ws.Range("J117").Value = 1
$ws.Range("K117").Value = "D"
$ws.Range("L117").Value = 8.5
$ws.Range("M117").Value = 4.75
$ws.Range("N117").Value = 1.363
$ws.Range("O117").Value = 8
$ws.Range("P117").Value = 4
$ws.Range("Q117").Value = 1.5
$ws.Range("R117").Value = 1
$ws.Range("S117").Value = 2.025
$ws.Range("T117").Value = 1.825
$ws.Range("U117").Value = 2.5
$ws.Range("V117").Value = 1.95
$ws.Range("W117").Value = 1.9
$ws.Range("X117").Value = -1
$ws.Range("Y117").Value = 3
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 1.025
$ws.Range("AB117").Value = -1
$ws.Range("AC117").Value = -1
$ws.Range("AD117").Value = 0.8999999999999999

# Row 118
$ws.Range("B118").Value = 6847102
$ws.Range("E118").Value = "Wolfsberger AC"
$ws.Range("F118").Value = "WSG Swarovski Tirol"
$ws.Range("G118").Value = 4
$ws.Range("H118").Value = 1
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = "H"
$ws.Range("L118").Value = 1.6
$ws.Range("M118").Value = 4.2
$ws.Range("N118").Value = 5.25
$ws.Range("O118").Value = 1.75
$ws.Range("P118").Value = 3.75
$ws.Range("Q118").Value = 4.5
$ws.Range("R118").Value = -0.75
$ws.Range("S118").Value = 2
$ws.Range("T118").Value = 1.85
$ws.Range("U118").Value = 2.5
$ws.Range("V118").Value = 1.85
$ws.Range("W118").Value = 2
$ws.Range("X118").Value = 0.75
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = -1
$ws.Range("AA118").Value = 1
$ws.Range("AB118").Value = -1
$ws.Range("AC118").Value = 0.8500000000000001
$ws.Range("AD118").Value = -1

# Row 159
$ws.Range("B159").Value = 7948273
$ws.Range("E159").Value = "WSG Swarovski Tirol"
$ws.Range("F159").Value = "FC Blau Weiss Linz"
$ws.Range("G159").Value = 2
$ws.Range("H159").Value = 1
$ws.Range("I159").Value = 1
$ws.Range("J159").Value = 0
$ws.Range("K159").Value = "H"
$ws.Range("L159").Value = 2.9
$ws.Range("M159").Value = 3.3
$ws.Range("N159").Value = 2.4
$ws.Range("O159").Value = 2.55
$ws.Range("P159").Value = 3.1
$ws.Range("Q159").Value = 3
$ws.Range("R159").Value = 0
$ws.Range("S159").Value = 1.775
$ws.Range("T159").Value = 2.1
$ws.Range("U159").Value = 2.25
$ws.Range("V159").Value = 2.1
$ws.Range("W159").Value = 1.775
$ws.Range("X159").Value = 1.55
$ws.Range("Y159").Value = -1
$ws.Range("Z159").Value = -1
$ws.Range("AA159").Value = 0.7749999999999999
$ws.Range("AB159").Value = -1
$ws.Range("AC159").Value = 1.1
$ws.Range("AD159").Value = -1

# Row 160
$ws.Range("B160").Value = 7948271
$ws.Range("E160").Value = "Austria Lustenau"
$ws.Range("F160").Value = "SCR Altach"
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 1
$ws.Range("I160").Value = 0
$ws.Range("J160").Value = 0
$ws.Range("K160").Value = "A"
$ws.Range("L160").Value = 3.75
$ws.Range("M160").Value = 3.4
$ws.Range("N160").Value = 2
$ws.Range("O160").Value = 5
$ws.Range("P160").Value = 3.4
$ws.Range("Q160").Value = 1.8
$ws.Range("R160").Value = 0.75
$ws.Range("S160").Value = 1.8
$ws.Range("T160").Value = 2.05
$ws.Range("U160").Value = 2
$ws.Range("V160").Value = 1.8
$ws.Range("W160").Value = 2.05
$ws.Range("X160").Value = -1
$ws.Range("Y160").Value = -1
$ws.Range("Z160").Value = 0.8
$ws.Range("AA160").Value = -0.5
$ws.Range("AB160").Value = 0.5249999999999999
$ws.Range("AC160").Value = -1
$ws.Range("AD160").Value = 1.05

# Row 165
$ws.Range("B165").Value = 7948276
$ws.Range("E165").Value = "SCR Altach"
$ws.Range("F165").Value = "FC Blau Weiss Linz"
$ws.Range("G165").Value = 2
$ws.Range("H165").Value = 2
$ws.Range("I165").Value = 1
$ws.Range("J165").Value = 1
$ws.Range("K165").Value = "D"
$ws.Range("L165").Value = 2.1
$ws.Range("M165").Value = 3.2
$ws.Range("N165").Value = 3.6
$ws.Range("O165").Value = 2
$ws.Range("P165").Value = 3
$ws.Range("Q165").Value = 4.75
$ws.Range("R165").Value = -0.5
$ws.Range("S165").Value = 2
$ws.Range("T165").Value = 1.85
$ws.Range("U165").Value = 1.75
$ws.Range("V165").Value = 1.775
$ws.Range("W165").Value = 2.1
$ws.Range("X165").Value = -1
$ws.Range("Y165").Value = 2
$ws.Range("Z165").Value = -1
$ws.Range("AA165").Value = -1
$ws.Range("AB165").Value = 0.8500000000000001
$ws.Range("AC165").Value = 0.7749999999999999
$ws.Range("AD165").Value = -1

# Row 166
$ws.Range("B166").Value = 7948275
$ws.Range("E166").Value = "FK Austria Vienna"
$ws.Range("F166").Value = "WSG Swarovski Tirol"
$ws.Range("G166").Value = 3
$ws.Range("H166").Value = 0
$ws.Range("I166").Value = 2
$ws.Range("J166").Value = 0
$ws.Range("K166").Value = "H"
$ws.Range("L166").Value = 1.65
$ws.Range("M166").Value = 3.8
$ws.Range("N166").Value = 5
$ws.Range("O166").Value = 1.8
$ws.Range("P166").Value = 3.6
$ws.Range("Q166").Value = 4.333
$ws.Range("R166").Value = -0.5
$ws.Range("S166").Value = 1.8
$ws.Range("T166").Value = 2.05
$ws.Range("U166").Value = 2.5
$ws.Range("V166").Value = 2.05
$ws.Range("W166").Value = 1.8
$ws.Range("X166").Value = 0.8
$ws.Range("Y166").Value = -1
$ws.Range("Z166").Value = -1
$ws.Range("AA166").Value = 0.8
$ws.Range("AB166").Value = -1
$ws.Range("AC166").Value = 1.05
$ws.Range("AD166").Value = -1

# Row 173
$ws.Range("B173").Value = 7948256
$ws.Range("E173").Value = "Hartberg"
$ws.Range("F173").Value = "LASK Linz"
$ws.Range("G173").Value = 1
$ws.Range("H173").Value = 2
$ws.Range("I173").Value = 1
$ws.Range("J173").Value = 1
$ws.Range("K173").Value = "A"
$ws.Range("L173").Value = 3
$ws.Range("M173").Value = 3.5
$ws.Range("N173").Value = 2.25
$ws.Range("O173").Value = 3.2
$ws.Range("P173").Value = 3.5
$ws.Range("Q173").Value = 2.2
$ws.Range("R173").Value = 0.25
$ws.Range("S173").Value = 1.975
$ws.Range("T173").Value = 1.875
$ws.Range("U173").Value = 2.25
$ws.Range("V173").Value = 1.8
$ws.Range("W173").Value = 2.05
$ws.Range("X173").Value = -1
$ws.Range("Y173").Value = -1
$ws.Range("Z173").Value = 1.2
$ws.Range("AA173").Value = -1
$ws.Range("AB173").Value = 0.875
$ws.Range("AC173").Value = 0.8
$ws.Range("AD173").Value = -1

# Row 174
$ws.Range("B174").Value = 7948255
$ws.Range("E174").Value = "Austria Klagenfurt"
$ws.Range("F174").Value = "FC Salzburg"
$ws.Range("G174").Value = 4
$ws.Range("H174").Value = 3
$ws.Range("I174").Value = 0
$ws.Range("J174").Value = 2
$ws.Range("K174").Value = "H"
$ws.Range("L174").Value = 6.5
$ws.Range("M174").Value = 4.5
$ws.Range("N174").Value = 1.444
$ws.Range("O174").Value = 9
$ws.Range("P174").Value = 5.25
$ws.Range("Q174").Value = 1.333
$ws.Range("R174").Value = 1.5
$ws.Range("S174").Value = 1.925
$ws.Range("T174").Value = 1.925
$ws.Range("U174").Value = 3
$ws.Range("V174").Value = 1.825
$ws.Range("W174").Value = 2.025
$ws.Range("X174").Value = 8
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = -1
$ws.Range("AA174").Value = 0.925
$ws.Range("AB174").Value = -1
$ws.Range("AC174").Value = 0.825
$ws.Range("AD174").Value = -1

# Row 183
$ws.Range("B183").Value = 7948285
$ws.Range("E183").Value = "SCR Altach"
$ws.Range("F183").Value = "Wolfsberger AC"
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 1
$ws.Range("I183").Value = 0
$ws.Range("J183").Value = 1
$ws.Range("K183").Value = "A"
$ws.Range("L183").Value = 2.15
$ws.Range("M183").Value = 3.2
$ws.Range("N183").Value = 3.5
$ws.Range("O183").Value = 2.1
$ws.Range("P183").Value = 3.25
$ws.Range("Q183").Value = 3.8
$ws.Range("R183").Value = -0.5
$ws.Range("S183").Value = 2.05
$ws.Range("T183").Value = 1.75
$ws.Range("U183").Value = 2.25
$ws.Range("V183").Value = 2.1
$ws.Range("W183").Value = 1.775
$ws.Range("X183").Value = -1
$ws.Range("Y183").Value = -1
$ws.Range("Z183").Value = 2.8
$ws.Range("AA183").Value = -1
$ws.Range("AB183").Value = 0.75
$ws.Range("AC183").Value = -1
$ws.Range("AD183").Value = 0.7749999999999999

# Row 184
$ws.Range("B184").Value = 7948284
$ws.Range("E184").Value = "FC Blau Weiss Linz"
$ws.Range("F184").Value = "WSG Swarovski Tirol"
$ws.Range("G184").Value = 3
$ws.Range("H184").Value = 2
$ws.Range("I184").Value = 2
$ws.Range("J184").Value = 1
$ws.Range("K184").Value = "H"
$ws.Range("L184").Value = 2.15
$ws.Range("M184").Value = 3.2
$ws.Range("N184").Value = 3.5
$ws.Range("O184").Value = 2.15
$ws.Range("P184").Value = 3.2
$ws.Range("Q184").Value = 3.75
$ws.Range("R184").Value = -0.25
$ws.Range("S184").Value = 1.825
$ws.Range("T184").Value = 2.025
$ws.Range("U184").Value = 2.25
$ws.Range("V184").Value = 1.95
$ws.Range("W184").Value = 1.9
$ws.Range("X184").Value = 1.15
$ws.Range("Y184").Value = -1
$ws.Range("Z184").Value = -1
$ws.Range("AA184").Value = 0.825
$ws.Range("AB184").Value = -1
$ws.Range("AC184").Value = 0.95
$ws.Range("AD184").Value = -1

# Row 185
$ws.Range("B185").Value = 7948259
$ws.Range("E185").Value = "Austria Klagenfurt"
$ws.Range("F185").Value = "LASK Linz"
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 2
$ws.Range("I185").Value = 0
$ws.Range("J185").Value = 1
$ws.Range("K185").Value = "A"
$ws.Range("L185").Value = 3.75
$ws.Range("M185").Value = 3.4
$ws.Range("N185").Value = 2
$ws.Range("O185").Value = 4.2
$ws.Range("P185").Value = 3.8
$ws.Range("Q185").Value = 1.8
$ws.Range("R185").Value = 0.75
$ws.Range("S185").Value = 1.825
$ws.Range("T185").Value = 2.025
$ws.Range("U185").Value = 3
$ws.Range("V185").Value = 2.025
$ws.Range("W185").Value = 1.825
$ws.Range("X185").Value = -1
$ws.Range("Y185").Value = -1
$ws.Range("Z185").Value = 0.8
$ws.Range("AA185").Value = -1
$ws.Range("AB185").Value = 1.025
$ws.Range("AC185").Value = -1
$ws.Range("AD185").Value = 0.825

# Row 186
$ws.Range("B186").Value = 7947239
$ws.Range("E186").Value = "SK Sturm Graz"
$ws.Range("F186").Value = "Hartberg"
$ws.Range("G186").Value = 1
$ws.Range("H186").Value = 1
$ws.Range("I186").Value = 0
$ws.Range("J186").Value = 1
$ws.Range("K186").Value = "D"
$ws.Range("L186").Value = 1.444
$ws.Range("M186").Value = 4.2
$ws.Range("N186").Value = 7.5
$ws.Range("O186").Value = 1.45
$ws.Range("P186").Value = 4.5
$ws.Range("Q186").Value = 7
$ws.Range("R186").Value = -1.25
$ws.Range("S186").Value = 2.05
$ws.Range("T186").Value = 1.8
$ws.Range("U186").Value = 3
$ws.Range("V186").Value = 2
$ws.Range("W186").Value = 1.85
$ws.Range("X186").Value = -1
$ws.Range("Y186").Value = 3.5
$ws.Range("Z186").Value = -1
$ws.Range("AA186").Value = -1
$ws.Range("AB186").Value = 0.8
$ws.Range("AC186").Value = -1
$ws.Range("AD186").Value = 0.8500000000000001

# Row 195
$ws.Range("B195").Value = 7948290
$ws.Range("E195").Value = "FC Blau Weiss Linz"
$ws.Range("F195").Value = "FK Austria Vienna"
$ws.Range("G195").Value = 1
$ws.Range("H195").Value = 2
$ws.Range("I195").Value = 0
$ws.Range("J195").Value = 1
$ws.Range("K195").Value = "A"
$ws.Range("L195").Value = 3.1
$ws.Range("M195").Value = 3.5
$ws.Range("N195").Value = 2.2
$ws.Range("O195").Value = 3.25
$ws.Range("P195").Value = 3.8
$ws.Range("Q195").Value = 2.05
$ws.Range("R195").Value = 0.5
$ws.Range("S195").Value = 1.825
$ws.Range("T195").Value = 2.025
$ws.Range("U195").Value = 2.75
$ws.Range("V195").Value = 1.975
$ws.Range("W195").Value = 1.875
$ws.Range("X195").Value = -1
$ws.Range("Y195").Value = -1
$ws.Range("Z195").Value = 1.05
$ws.Range("AA195").Value = -1
$ws.Range("AB195").Value = 1.025
$ws.Range("AC195").Value = 0.4875
$ws.Range("AD195").Value = -0.5

# Row 196
$ws.Range("B196").Value = 7948291
$ws.Range("E196").Value = "SCR Altach"
$ws.Range("F196").Value = "Austria Lustenau"
$ws.Range("G196").Value = 2
$ws.Range("H196").Value = 2
$ws.Range("I196").Value = 1
$ws.Range("J196").Value = 1
$ws.Range("K196").Value = "D"
$ws.Range("L196").Value = 1.666
$ws.Range("M196").Value = 4
$ws.Range("N196").Value = 4.75
$ws.Range("O196").Value = 1.6
$ws.Range("P196").Value = 4.2
$ws.Range("Q196").Value = 5
$ws.Range("R196").Value = -1
$ws.Range("S196").Value = 2.05
$ws.Range("T196").Value = 1.8
$ws.Range("U196").Value = 2.5
$ws.Range("V196").Value = 1.875
$ws.Range("W196").Value = 1.975
$ws.Range("X196").Value = -1
$ws.Range("Y196").Value = 3.2
$ws.Range("Z196").Value = -1
$ws.Range("AA196").Value = -1
$ws.Range("AB196").Value = 0.8
$ws.Range("AC196").Value = 0.875
$ws.Range("AD196").Value = -1
